# Append a new "S7" Activity block (with its Exchanges sub-table) to Sheet1,
# mirroring the existing "Test" Activity block (rows 69-79) but for the new
# S7 activity whose single exchange references the "Landfill_saint sophie"
# technosphere flow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Activity header block (rows 84-88) ---
$ws.Range("A84").Value = "Activity"
$ws.Range("B84").Value = "S7"

$ws.Range("A85").Value = "comment"
$ws.Range("B85").Value = "S7"

$ws.Range("A86").Value = "location"
$ws.Range("B86").Value = "CA-QC"

$ws.Range("A87").Value = "production amount"
$ws.Range("B87").NumberFormat = "@"
$ws.Range("B87").Value = "1"

$ws.Range("A88").Value = "unit"
$ws.Range("B88").Value = "tonne"

# row 89 left blank (separator row)

# --- Exchanges sub-table (rows 90-93) ---
$ws.Range("A90").Value = "Exchanges"

$ws.Range("A91").Value = "name"
$ws.Range("B91").Value = "reference product"
$ws.Range("C91").Value = "unit"
$ws.Range("D91").Value = "amount"
$ws.Range("E91").Value = "location"
$ws.Range("F91").Value = "database"
$ws.Range("G91").Value = "type"
$ws.Range("H91").Value = "categories"
$ws.Range("I91").Value = "comment"

# Production exchange (self-reference)
$ws.Range("A92").Value = "S7"
$ws.Range("B92").Value = "OFMSW"
$ws.Range("C92").Value = "tonne"
$ws.Range("D92").NumberFormat = "@"
$ws.Range("D92").Value = "1"
$ws.Range("E92").Value = "CA-QC"
$ws.Range("F92").Value = "Scenarios"
$ws.Range("G92").Value = "production"

# Technosphere exchange
$ws.Range("A93").Value = "Landfill_saint sophie"
$ws.Range("B93").Value = "OFMSW"
$ws.Range("C93").Value = "tonne"
$ws.Range("D93").Value = 1
$ws.Range("E93").Value = "CA-QC"
$ws.Range("F93").Value = "OWM Facilities"
$ws.Range("G93").Value = "technosphere"
